# Proyecto trimestral: interfaz grafica avanzada, registro de usuarios,
# productos, edicion de usuarios, modulo de caja y asignacion de base.
#
# Rebuilds the "Inventario" sheet: new headers/columns (disponibilidad,
# categoria) and a 3-article product table (Arbeja/Frijol/Garbanzo).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inventario")

function Set-TextCell($addr, $val) {
    # Force number-format "Text" before writing so numeric-looking values
    # ("1000", "0", "2022-10-25", ...) are stored as shared strings (t="s"),
    # matching the source workbook instead of being auto-coerced to
    # numbers/dates.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

# ---- Row 1: headers ------------------------------------------------------
$ws.Range("A1").Value = "cod_articulo"
$ws.Range("B1").Value = "nombre_producto"
$ws.Range("C1").Value = "disponibilidad"
$ws.Range("D1").Value = "valor_unitario"
$ws.Range("E1").Value = "iva_producto"
$ws.Range("F1").Value = "categoria"
$ws.Range("G1").Value = "vencimiento_producto"

# ---- Row 2: Arbeja (A2 written last so shared-string order matches) -----
Set-TextCell "B2" "Arbeja"
Set-TextCell "C2" "1000"
Set-TextCell "D2" "100"
Set-TextCell "E2" "0"
Set-TextCell "F2" "1"
Set-TextCell "G2" "2022-10-25"
Set-TextCell "A2" "0001"

# ---- Row 3: Frijol --------------------------------------------------------
Set-TextCell "A3" "0002"
Set-TextCell "B3" "Frijol"
Set-TextCell "C3" "2000"
Set-TextCell "D3" "50"
Set-TextCell "E3" "0"
Set-TextCell "F3" "1"
Set-TextCell "G3" "2023-11-22"

# ---- Row 4: Garbanzo -------------------------------------------------------
Set-TextCell "A4" "0003"
Set-TextCell "B4" "Garbanzo"
Set-TextCell "C4" "4000"
Set-TextCell "D4" "40"
Set-TextCell "E4" "0"
Set-TextCell "F4" "1"
Set-TextCell "G4" "2023-11-22"

# ---- Formatting: table body (A1:G4) --------------------------------------
$body = $ws.Range("A1:G4")
$body.NumberFormat = "@"
$body.HorizontalAlignment = 1
$body.VerticalAlignment = -4107

# ---- Formatting: decorative empty cells (I1:I3) --------------------------
$deco = $ws.Range("I1:I3")
$deco.HorizontalAlignment = -4108
$deco.VerticalAlignment = -4108
$deco.WrapText = $true

# ---- Column widths ---------------------------------------------------------
$ws.Columns("E").ColumnWidth = 20
$ws.Columns("G").ColumnWidth = 21

# ---- Row heights ------------------------------------------------------------
$ws.Rows(1).RowHeight = 26.25

# ---- Selection --------------------------------------------------------------
$ws.Range("G4").Select()
